$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 848.65
$ws.Range("I19").Value = 716.25
$ws.Range("J19").Value = 936.9167
$ws.Range("K19").Value = 716.25
$ws.Range("L19").Value = 936.9167
$ws.Range("M19").Value = -541.25
$ws.Range("N19").Value = -1286.9167
$ws.Range("H33").Value = 410.08334
$ws.Range("I33").Value = 389.7619
$ws.Range("K33").Value = 389.7619
$ws.Range("M33").Value = -160.7619
$ws.Range("H70").Value = 1875.9286
$ws.Range("I70").Value = 1852
$ws.Range("J70").Value = 1889.2222
$ws.Range("K70").Value = 5556
$ws.Range("L70").Value = 5667.6666
$ws.Range("M70").Value = -5286
$ws.Range("N70").Value = -6207.6666
$ws.Range("H73").Value = 1875.9286
$ws.Range("I73").Value = 1852
$ws.Range("J73").Value = 1889.2222
$ws.Range("K73").Value = 5556
$ws.Range("L73").Value = 5667.6666
$ws.Range("M73").Value = -4620
$ws.Range("N73").Value = -7539.6666
$ws.Range("H98").Value = 2022.4822
$ws.Range("I98").Value = 2125.9795
$ws.Range("J98").Value = 1298
$ws.Range("K98").Value = 2125.9795
$ws.Range("L98").Value = 1298
$ws.Range("M98").Value = -627.9794999999999
$ws.Range("N98").Value = -4294
$ws.Range("H107").Value = 4363.125
$ws.Range("I107").Value = 6968.3335
$ws.Range("J107").Value = 2800
$ws.Range("K107").Value = 6968.3335
$ws.Range("L107").Value = 2800
$ws.Range("M107").Value = -5048.3335
$ws.Range("N107").Value = -6640
$ws.Range("H116").Value = 3351.875
$ws.Range("J116").Value = 3822.5
$ws.Range("L116").Value = 3822.5
$ws.Range("N116").Value = -10706.5
$ws.Range("H122").Value = 2022.4822
$ws.Range("I122").Value = 2125.9795
$ws.Range("J122").Value = 1298
$ws.Range("K122").Value = 6377.9385
$ws.Range("L122").Value = 3894
$ws.Range("M122").Value = -3927.9385
$ws.Range("N122").Value = -8794
$ws.Range("H135").Value = 3862.2222
$ws.Range("I135").Value = 756
$ws.Range("J135").Value = 6347.2
$ws.Range("K135").Value = 6804
$ws.Range("L135").Value = 57124.8
$ws.Range("M135").Value = -4269
$ws.Range("N135").Value = -62194.8

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 23100
$ws.Range("J24").Value = 23100
$ws.Range("L24").Value = 23100
$ws.Range("N24").Value = -23848
$ws.Range("H61").Value = 142859580
$ws.Range("I61").Value = 166669010
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 166669010
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -166668798
$ws.Range("N61").Value = -3424
$ws.Range("H74").Value = 1697.9286
$ws.Range("I74").Value = 1197.3636
$ws.Range("K74").Value = 1197.3636
$ws.Range("M74").Value = -323.3635999999999
$ws.Range("H77").Value = 1697.9286
$ws.Range("I77").Value = 1197.3636
$ws.Range("K77").Value = 5986.817999999999
$ws.Range("M77").Value = -1618.817999999999
$ws.Range("H100").Value = 23100
$ws.Range("J100").Value = 23100
$ws.Range("L100").Value = 23100
$ws.Range("N100").Value = -25264
$ws.Range("H136").Value = 142859580
$ws.Range("I136").Value = 166669010
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 500007030
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -500004480
$ws.Range("N136").Value = -14100

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 142858510
$ws.Range("I105").Value = 166667760
$ws.Range("J105").Value = 3011
$ws.Range("K105").Value = 166667760
$ws.Range("L105").Value = 3011
$ws.Range("M105").Value = -166666013
$ws.Range("N105").Value = -6505
$ws.Range("H134").Value = 12656.1
$ws.Range("I134").Value = 1937.2858
$ws.Range("J134").Value = 37666.668
$ws.Range("K134").Value = 5811.857400000001
$ws.Range("L134").Value = 113000.004
$ws.Range("M134").Value = -3276.857400000001
$ws.Range("N134").Value = -118070.004

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1294.1621
$ws.Range("I31").Value = 1093.9706
$ws.Range("K31").Value = 1093.9706
$ws.Range("M31").Value = -798.9706000000001
$ws.Range("H34").Value = 1294.1621
$ws.Range("I34").Value = 1093.9706
$ws.Range("K34").Value = 1093.9706
$ws.Range("M34").Value = -891.9706000000001
$ws.Range("H100").Value = 47940
$ws.Range("J100").Value = 47940
$ws.Range("L100").Value = 47940
$ws.Range("N100").Value = -50104
$ws.Range("H138").Value = 172345
$ws.Range("J138").Value = 172345
$ws.Range("L138").Value = 172345
$ws.Range("N138").Value = -182625

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 22730920
$ws.Range("J131").Value = 3983.475
$ws.Range("L131").Value = 11950.425
$ws.Range("N131").Value = -22030.425

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 37503170
$ws.Range("I70").Value = 27781056
$ws.Range("J70").Value = 66669500
$ws.Range("K70").Value = 27781056
$ws.Range("L70").Value = 66669500
$ws.Range("M70").Value = -27780786
$ws.Range("N70").Value = -66670040
$ws.Range("H73").Value = 37503170
$ws.Range("I73").Value = 27781056
$ws.Range("J73").Value = 66669500
$ws.Range("K73").Value = 27781056
$ws.Range("L73").Value = 66669500
$ws.Range("M73").Value = -27780120
$ws.Range("N73").Value = -66671372
$ws.Range("H80").Value = 2747.0667
$ws.Range("I80").Value = 1657.1428
$ws.Range("K80").Value = 1657.1428
$ws.Range("M80").Value = -659.1428000000001
$ws.Range("H83").Value = 2747.0667
$ws.Range("I83").Value = 1657.1428
$ws.Range("K83").Value = 8285.714
$ws.Range("M83").Value = -3293.714
$ws.Range("H102").Value = 1874.75
$ws.Range("I102").Value = 1333
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 1333
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = 289
$ws.Range("N102").Value = -6744
$ws.Range("H113").Value = 1867.3334
$ws.Range("I113").Value = 1867.3334
$ws.Range("K113").Value = 1867.3334
$ws.Range("M113").Value = 302.6666

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1901.3636
$ws.Range("I7").Value = 1750
$ws.Range("J7").Value = 2305
$ws.Range("K7").Value = 1750
$ws.Range("L7").Value = 2305
$ws.Range("M7").Value = -1638
$ws.Range("N7").Value = -2529
$ws.Range("H101").Value = 6833
$ws.Range("J101").Value = 6833
$ws.Range("L101").Value = 6833
$ws.Range("N101").Value = -13323
$ws.Range("H122").Value = 16677967
$ws.Range("I122").Value = 20843642
$ws.Range("J122").Value = 15269.667
$ws.Range("K122").Value = 62530926
$ws.Range("L122").Value = 45809.001
$ws.Range("M122").Value = -62528476
$ws.Range("N122").Value = -50709.001
$ws.Range("H126").Value = 1901.3636
$ws.Range("I126").Value = 1750
$ws.Range("J126").Value = 2305
$ws.Range("K126").Value = 5250
$ws.Range("L126").Value = 6915
$ws.Range("M126").Value = -2780
$ws.Range("N126").Value = -11855
$ws.Range("H136").Value = 21800.6
$ws.Range("I136").Value = 26500.75
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 79502.25
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -76952.25
$ws.Range("N136").Value = -14100

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 27783222
$ws.Range("J62").Value = 20003
$ws.Range("L62").Value = 20003
$ws.Range("N62").Value = -21251
$ws.Range("H65").Value = 27783222
$ws.Range("J65").Value = 20003
$ws.Range("L65").Value = 100015
$ws.Range("N65").Value = -106255
$ws.Range("H96").Value = 2915.4614
$ws.Range("I96").Value = 1866.3334
$ws.Range("K96").Value = 1866.3334
$ws.Range("M96").Value = -493.3334
$ws.Range("H113").Value = 462.8
$ws.Range("I113").Value = 291.33334
$ws.Range("J113").Value = 720
$ws.Range("K113").Value = 874.0000200000001
$ws.Range("L113").Value = 2160
$ws.Range("M113").Value = 1295.99998
$ws.Range("N113").Value = -6500
$ws.Range("H132").Value = 2352.7932
$ws.Range("I132").Value = 1266.3572
$ws.Range("J132").Value = 3366.8
$ws.Range("K132").Value = 3799.0716
$ws.Range("L132").Value = 10100.4
$ws.Range("M132").Value = -1269.0716
$ws.Range("N132").Value = -15160.4
$ws.Range("H136").Value = 1098.037
$ws.Range("I136").Value = 1118.6471
$ws.Range("J136").Value = 1063
$ws.Range("K136").Value = 3355.9413
$ws.Range("L136").Value = 3189
$ws.Range("M136").Value = -805.9412999999995
$ws.Range("N136").Value = -8289
